# ---------------------------------------------------------------------------
# Applies the "Built site for gh-pages" commit:
#   1. Insert a new "Author" paragraph ("Dr. Lennart Wittkuhn") as the very
#      first paragraph of the document body.
#   2. Switch the timetable table to a fixed column layout
#      (<w:tblLayout w:type="fixed"/>).
#   3. Add a new "Abstract Title" paragraph style (custom style id
#      "AbstractTitle"), based on Normal, followed by Abstract.
#   4. Reduce the "Abstract" style's space-before from 15pt (300 twips) to
#      5pt (100 twips).
#   5. Add a new "Footnote Block Text" paragraph style (style id
#      "FootnoteBlockText"), based on the Footnote Text style.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. New Author paragraph at the very top of the document, built from three
#    separate runs ("Dr. Lennart" / " " / "Wittkuhn") exactly like pandoc
#    would emit for the "Dr. Lennart Wittkuhn" author metadata entry.
# ---------------------------------------------------------------------------
$authorXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="Author"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Dr. Lennart</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">Wittkuhn</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$startRange = $d.Range(0, 0)
$startRange.InsertXML($authorXml)

# ---------------------------------------------------------------------------
# 2. Fixed table layout for the (single) timetable table.
# ---------------------------------------------------------------------------
$table = $d.Tables(1)
$table.AllowAutoFit = $false

# ---------------------------------------------------------------------------
# 3. New "Abstract Title" style, inserted ahead of the "Abstract" style.
# ---------------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# ---------------------------------------------------------------------------
# 4. "Abstract" style: space-before 300 twips -> 100 twips (15pt -> 5pt).
# ---------------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------------
# 5. New "Footnote Block Text" style, based on Footnote Text.
# ---------------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "done"
